# Apply SF calibration-result updates:
#  - new "note" column (L) with two annotations
#  - two new calibration rows (3 and 4)
#  - widen columns for the new "simulation_roll_outs" (F) and "supermarket" (K) data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header for the new "note" column ---
$ws.Range("L1").Value = "note"

# --- Row 3: new calibration run ---
# dates: reuse the existing formatted cells' look by copying their formats first
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> custom yyyy-mm-dd;@ style
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> m/d/yyyy style

$ws.Range("A3").Value = 43894
$ws.Range("B3").Value = 43909
$ws.Range("C3").Value = 100
$ws.Range("D3").Value = 20
$ws.Range("E3").Value = 100
$ws.Range("F3").Value = 40
$ws.Range("G3").Value = 0.0151399550959467
$ws.Range("H3").Value = 0.111006572842597
$ws.Range("I3").Value = 0.0043372111395001403
$ws.Range("J3").Value = 0.23195971548557201
$ws.Range("K3").Value = 1.35917520523071
$ws.Range("L3").Value = "subdivide offices"

# --- Row 4: new calibration run (sparse - several columns left blank) ---
$ws.Range("B2").Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> m/d/yyyy style
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> m/d/yyyy style

$ws.Range("A4").Value = 43894
$ws.Range("B4").Value = 43945
$ws.Range("C4").Value = 100
$ws.Range("D4").Value = 20
$ws.Range("F4").Value = 40
$ws.Range("H4").Value = 0.04
$ws.Range("I4").Value = 0.1
$ws.Range("J4").Value = 0.04
$ws.Range("K4").Value = "0.2?"

# --- Column widths for the newly-used columns ---
$ws.Range("F1").Value = "simulation_roll_outs"
$ws.Columns("F:F").ColumnWidth = 16
$ws.Columns("K:K").ColumnWidth = 10.76
